$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "CHUYE" -> "CHUYEN" in A9 (shared string used by that cell)
$ws.Range("A9").Value = "CHUYEN"

# A9:B9 currently carry an explicit "apply fill" formatting flag that makes
# their style distinct from the rest of the table (A2:B8). Clearing the
# explicit fill (it was already "no fill") realigns them onto the same
# style used elsewhere in the table.
$ws.Range("A9:B9").Interior.Pattern = -4142

# Move the active selection to B14 (below the used range), matching the
# cursor position left after the edit.
$ws.Range("B14").Select() | Out-Null
